$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert two rows at the top, pushing existing data down.
$ws.Rows("1:2").Insert()

# The table (originally A1:C18) needs to move down with the data to A3:C20.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:C20"))

# Title row
$ws.Range("A1").Value = "DICTIONARY OF STORED PROCEDURES"
$ws.Range("A1").Font.Bold = $true

# Rename the sheet
$ws.Name = "Dic SP"

# Set print area
$ws.PageSetup.PrintArea = "`$A`$1:`$C`$20"

# Leave the cursor where the author left it when they saved.
[void]$ws.Range("B22").Select()
